# fix(import): add antenne column
# Adds a new "antenne" header (column K) with a "MONTREUIL" sample value
# in row 2, matching the rest of the import-template headers/data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "antenne"
$ws.Range("K2").Value = "MONTREUIL"

$ws.Range("K3").Select()
